$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the sheet "My Series" -> "Data"
$ws.Name = "Data"

# 2. Update the CEIC metadata comment payload attached to cell A1
$comment = $ws.Range("A1").Comment
$comment.Text("g0AAAB+LCAAAAAAAAAPtnOtvG1UWwP+VUaSVQKozM3biONbtID/y8DZOQuw2LV+q8cxNMpvxjHceSf2t1QpaWApiKSBe210QsIsoBakLpQnd/wXVbvqp/8Kee+887XHqScuKIqMKMudxH+fee85vbjOgly61dW4PW7ZmGqenxGlhisOGYqqasX16ynW2MmJ+6iUJLVxSsL4uW3IbO2DMgZdhFy/Z2umpHcfpFHl+f39/ej83bVrbfFYQRP58faWh7OC2nNEM25ENBU8FXuqTvaYkVFHbdezIquzIzPP0VK1Rm65gTamCrC4b8ja2psuurRnYthcMR3M0bBNPC8sOrlTr59jEpOx0flpE/JA8tCy7mq4yu5glk3t20C1uam0sZQWxkBHFTHa+KcwVhUJRzE4XxNwrvmNgiFZk22lga09TqKDhyO0OdRfF7LwwJxTELOITjaCtMAASWtPVDbyn2VitYF23U0WE9xawpDgw63TBFBAf8fUaOvkQliy5s9PUHB2ndV80LaxAoE7U9yreX7O8+DU7K6Bt7miW063K3dRtnbWxtdYh0UjnKqGqaTglHVvO2Q4sKlZhzUEhOZaLET9CGTpVNVuBnzXDxaq0Jet21CmmRJumtWt3ZAWvwoHlSRv7hm7KKuwsR7MdTQk7HVKgdcvsQIvQednU1UVo1TNOUAQt1wwIMem2bJq74eiSlIjuAbobYE3bsuObD8lRY8fcXzP0bsNt2YqltbBaLfvWiTpETp7nXXFtx2zDKEIRYrKIRFX5ep3vwj9w3gaVqIoVrS3r6zqE0pZy0FZMgEquY25pTsXU3bZh+yMbkKJNmFcTXwrmGTyjNVhig4TeNGqGb8+CnaiKO2yY+0Gfwwoaioi4ZCv+og8rBo2rIPMXcVhD14XMclHToRxEVyQije+Nxg7GTuLGYBpEMt8iKTBSubvqtltwyFpw0vZorzbiQz2C3Qo7HsYlCVA0MvRPUxCK9A+MI1CjBUMdbecrEXQX6UsiKXlAhGBOelmXjV2QbmrOzmrJn0uCBrEIjLQf1iE4vx1d7lJxEKWoDNUMRXdVzNJCzdiiW5SMjS3qSDUaEq3ASZeQbHSb3Q6kYVsrOvDD6SkozEXbsaD0T0mK6RqO1SX5A/Ge6ZN8bLdl0A5kfWyfLQv/2QXi6C66hlIx1fF7U1l0zhqaM/4ITddiSXF8Fxo9kh9du4pJpqGpf2x/Jc2cbCuVedvAbdPQlPGjDUEmo1dPMBHbP1Vje2B2vsa216G6s+JHzvrYbhbgIpS7VN2UbNtUNLpZveOhRvz5EUemirdkVwdUc6DQbge5d1CMSvbuoE1UhM5aup8BJQLCNpCworanFUAIQnvTitkmAh4AdLOB+Kg94SAFLxjbK7Kx7QJpBHllUB7kX1Ilm5Zs2GQ6AVgMpOJkI+TnKQY8Ekteay7dCCx5maBF/IAdauJ2x7RkvQ6B0Ra9bedREwBJXXZ2vCeobTpW/CDzoWvgFR+ZP/AnmdEixaZBDryXJgeE1IjMhTF3aBPKEJllHY6lXpF1rWWxrOqX8iQdLFiIiH7+JZNLiYv+GsB7F1TfM7hLYDx88OR0y4q+gm1gkkilxsZMITszn5sHsiHPiM54yTJtm6uabUxwjwOmU13F4V5Yqq6/CGWPBsUrPGPYRy3Rop/LaXcvw5jikrgBMMW2BqVi2DDQhA7Syy5kH2zp3Ygtm9GKqYBp7+Ofe6+9+vDGzQf3vupfPuhdPnx8eA3G+PjwdTYrZoaackvHdDTNcqEg5GZgIwUiRKLHUwImk6SyCxco+AbPyHtRow+VhVplaaVME0Yg9N1ZzeDJO2DXdMPHBhs+7YiuGe8vNTORmn4C8p5j2kgNksg72R6OW0f1oxxZLB4efPXw4NZIby9gIUyJ8/PzGSH3RNaCt1sxI2ZHsNZKLMkT45mMMJvJZiPGAzZog2X4IE41VcqJwrwoAtUFyVoN9m2S0aDKa6kpb/MDfkxUYfwTbIHos6+ku74J5yJQs3MQeWBB7P/nr0e334tZedH1JPFWYHCUU0hnvP9Am17daHKNtbMblQWuudAg+yTURexY48cYe70HJym6aYKTdoqDAg4Fi5uqy9Yp7o+ugU9xDdw5xcFb0FTkHMa2XJKUdfdUDQ+OGNKS22F5JuIQShMsw4ST5JGQdKiOxjYp+4TaBA8vKd36PMnBm0s1BFhpo17mWkawTakMxTRMFNF7XVy/8uDefx/cu9f/7p3e3SuxFrx+gssA2PZwuKKPwSmADOjVlwEJ2mzQeO4KFyP1xBOSF6d1UzMcW5pl70zeEwJXkbRG/4tqbShxtGEaL5APSNCybC9ccrxzLq0iPi6AcXZkqK5m+K4ZCFhKD+P66NO/9z++03//u6OrX/fe+FfvzfcfHtw8+vaf7BD23/uuf/1bL+kP1gU6FvIGy6CPo7ciCkcOJ0dqNffL5Xc5w3Q4QAzOpQnql8sfRhojA6UwErYMCBcMJD6EIdOoM/HjIkMJxhDzC1xYwa+QipYLLLyaZnY0JezklQxpipw+qnih1sy4NuZMoCeo5wPGofO4fp4Lq7Drc0KWvEyHInLPILdkOxL6Jd1syTrnK+iFw4BJzOt4h9CW9re0slYurYQmbBBrlootsg3ZD8hHSFJharb/5G+1iAS0AHqKq5M7oiGzYVXQciST8d5ly1ZJJUkw+WoiZoEqrmUxODK8q/qG2wH49S/mRuvpZWWEd1cZm0YJOHyuVeN6eI5ooS7G1URA9TQ1eSqWpmo2uddh+LpKQhM+gi52wQnh8G7jGXjtabAkPMk7C5ZlWonJJ9T4ZnUgZ8gofBjxwIauKaNsNVwrX+AnvGdD27OiMBej7Uq5yJXUPfLXH3aRW4eJstuRYc4eZXkMYS+PS9hDhkmE3cBtTTYMV9ZHQnal/Pjwo6M7P/a/ud87+KH/yc3etX8/fOsqCB/c/Rpq0IS0nxVpFwCz8+OQdu7/QtrC7PyTSdszmpD2uKS9LOtbmS6WgSsDIgYYPilfn6C50VS9nJqqEz2OpeoRGWccsL7+Ru+nO2Oy9fqZDa6dyNa+ZgRbP/zw/oO7V3s/3undvNd//3bv0y971z/r377xlJAtjoZs8WKkjAxDdn4C2RPIfr4gW5xA9vMC2QPJ5zcL2fAWUIhB9gvVyovLJhyaM+RfG9iRNZ0rm4ZKQNoCWCxydU0tcstnuBXN2OVKRU6YoygkCPNFLjedF/6QhOTPpt1jAL46OybADxkmAXwVhjeK3B8fXutd+aR/6zOWtR8fvv7oyw/7d28/+viH3o1ve5fv9a7dJRh/8GP/H4eU5289+uAOvUX/yLN8937/rS9K8BzOER4iwZtw/9NyP8R0NiMApc88ifuFebIGQv5X5/7CONxfmHB/Gu6nB/UUNzXLmeGvmJwA+NO0M5r0q6lJP9HjWNIfTE1jIH7/gy/G5PvlM9UkuKfiEWQP+az36l+ekuOzozk+ezFSqYY5XhTyhQnJT0j+uSL57ITknxeSH0g/v1mSnxXIXxtGSH5zvUbRGtJYlytZkE10chu+aJqA2Yvnitw5vF3kliyMDW4dy3YStadv4xhC38yPSehDhkmEvonx7khE77/+U+/tdxiB99+8CmkPoB0mA6xOgfzg6Msrvfde6717hTD55zfZD4vnyJX8jW+O3v6EiD/929H3r01w/BniuCiMcw1P7uAz4syv/QsvudwYv/DCjCY4Pi6Os1PJtTCcV4NUbNnhGrA3LFXupofy9K2NRvPN1Gie6HEsmg/lpDHY/NHl78dkczhCM+SqYPa0SOv0EKUPGIzg9dHNpMf23Ghsz12MlKUEbJ+b3L9PqP35ovbc74Davb5/99g+kH5+FWz3fmDfM3gzNKtYx07qDzJ977q5d2JfWPy0rjV7TVe9YKb7hfogLGED0a9SyUZJO5iSZQFGkQ/YUn9GSv7KoJZu9tKsWGgJs/kZJZ+dz83JhYI6l58rbAkFVVRzggogSRsl31eQJtrA2uTrhHSdADzG3cMvSTZkYztlayzA1JF8bQGRInEmW1izbOc8qZPeT0xyIZBcYOx8XpphgHyePV+QcnPzTAIWfLR5PjZOP6U47JNlU1/R2lrKbzAEP+3EG4GF73QYXaZdQih5q/gSsG+kBUjWrT9BOWMfLaVpjZ0jyPGBP/lQ0Na2d5y0A1MUJdsSlGxmJofzmRlZnsvIWNyCF6DcXF6ZmZvJ5WfJZ4Ze45DQNLyfshPeX7DwfyUg/Q8dqoh1g0AAAA==")

# 3. Custom number format (numFmtId 166) "0.000" -> "###0.000"
#    This range is exactly every cell currently carrying that style.
$ws.Range("B27:E380").NumberFormat = "###0.000"

# 4. D1: series display name dropped "SAR (China)"
$ws.Range("D1").Value = "(DC)Hong Kong Retail Bonds: Price: Mid: HK Link A: 07-05-2009: 3.60%"

# 5. A11 label text
$ws.Range("A11").Value = "Function Information"

# 6. D14: corrected Last Update Time value
$ws.Range("D14").Value = 41781

# 7. D20: Skewness value precision fix
$ws.Range("D20").Value = -1.49587483823416

# 8. C21: Kurtosis value precision fix
$ws.Range("C21").Value = 0.3006970042568957
